# Auto-generated Excel COM-interop script
# Fix: reformat floating point "Importe" amounts from ar-AR style (1.234,56) to
# invariant style (1234.56), and normalize commas to periods in a few
# "Razon social" name fields that were mis-scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name field fixes (commas -> periods) ---
$nameCells = @(
    @("E50", "FERNANDEZ. MARIO HUGO"),
    @("E55", "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"),
    @("E149", "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"),
    @("E82", "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"),
    @("E95", "BOFFELLI. MARIA INES"),
    @("E110", "RICCOTTI. MARIANA EDITH"),
    @("E161", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
)

foreach ($pair in $nameCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Amount field fixes (es-AR decimal formatting -> invariant en-US) ---
# These text values look like numbers, so Excel would silently coerce a plain
# ".Value = " assignment into a real number (losing trailing zeros, e.g. 7750.00
# -> 7750). Force the cell to Text format first, assign the literal string, then
# restore the cell style so no stray number-format is left behind.
$amountCells = @(
    @("H2", "7750.00"),
    @("H3", "1900.00"),
    @("H4", "4300.00"),
    @("H5", "5300.00"),
    @("H6", "1277283.00"),
    @("H7", "80000.00"),
    @("H8", "189120.00"),
    @("H9", "235008.00"),
    @("H10", "26.90"),
    @("H11", "4176.00"),
    @("H12", "3800.00"),
    @("H13", "3690.00"),
    @("H14", "100.00"),
    @("H15", "500.00"),
    @("H16", "8177.10"),
    @("H17", "8031.10"),
    @("H18", "1050.00"),
    @("H19", "5843.70"),
    @("H20", "402902.56"),
    @("H21", "70447.60"),
    @("H22", "89703.72"),
    @("H23", "32235.00"),
    @("H24", "3254.00"),
    @("H25", "24408.00"),
    @("H26", "47750.63"),
    @("H27", "645.00"),
    @("H28", "9802.72"),
    @("H29", "19999.46"),
    @("H30", "8260.00"),
    @("H31", "300.00"),
    @("H32", "1316.31"),
    @("H33", "3974.00"),
    @("H34", "217.20"),
    @("H35", "757.20"),
    @("H36", "1900.00"),
    @("H37", "2120.00"),
    @("H38", "80.00"),
    @("H39", "31663.82"),
    @("H40", "2849.58"),
    @("H41", "161.26"),
    @("H42", "58480.00"),
    @("H43", "676.40"),
    @("H44", "2839.47"),
    @("H45", "1541.38"),
    @("H46", "362.75"),
    @("H47", "4686.00"),
    @("H48", "7843.49"),
    @("H49", "2039.00"),
    @("H50", "670.00"),
    @("H51", "2178.00"),
    @("H52", "5750.00"),
    @("H53", "101899.14"),
    @("H54", "139292.25"),
    @("H55", "3416.00"),
    @("H56", "14460.84"),
    @("H57", "150.00"),
    @("H58", "964.80"),
    @("H59", "3477.00"),
    @("H60", "246.27"),
    @("H61", "103.38"),
    @("H62", "17332.63"),
    @("H63", "2239.17"),
    @("H64", "3640.00"),
    @("H65", "418.61"),
    @("H66", "2736.00"),
    @("H67", "830.00"),
    @("H68", "3850.00"),
    @("H69", "160.00"),
    @("H70", "564.00"),
    @("H71", "1060.00"),
    @("H72", "142.40"),
    @("H73", "4530.00"),
    @("H74", "2870.00"),
    @("H75", "100.00"),
    @("H76", "2400.00"),
    @("H77", "16393.00"),
    @("H78", "1460.00"),
    @("H79", "550.00"),
    @("H80", "2472.00"),
    @("H81", "20798.00"),
    @("H82", "6860.00"),
    @("H83", "510.00"),
    @("H84", "3346.78"),
    @("H85", "9147.24"),
    @("H86", "415230.93"),
    @("H87", "53261.49"),
    @("H88", "43595.60"),
    @("H89", "2930.00"),
    @("H90", "83.56"),
    @("H91", "1950.00"),
    @("H92", "43385.31"),
    @("H93", "273.00"),
    @("H94", "653.40"),
    @("H95", "63852.80"),
    @("H96", "42.00"),
    @("H97", "225.00"),
    @("H98", "132665.88"),
    @("H99", "1594.00"),
    @("H100", "228.00"),
    @("H101", "285.00"),
    @("H102", "1334.00"),
    @("H103", "6606.44"),
    @("H104", "3300.00"),
    @("H105", "12404.32"),
    @("H106", "12263.50"),
    @("H107", "1057.80"),
    @("H108", "4350.00"),
    @("H109", "2900.00"),
    @("H110", "2000.00"),
    @("H111", "11200.00"),
    @("H112", "796.00"),
    @("H113", "6829.39"),
    @("H114", "2000.00"),
    @("H115", "5112.43"),
    @("H116", "73430.00"),
    @("H117", "2400.00"),
    @("H118", "4000.00"),
    @("H119", "3000.00"),
    @("H120", "19000.00"),
    @("H121", "5000.00"),
    @("H122", "4808.14"),
    @("H123", "590.00"),
    @("H124", "2123.90"),
    @("H125", "8534.60"),
    @("H126", "7380.00"),
    @("H127", "347.29"),
    @("H128", "987.19"),
    @("H129", "1504.81"),
    @("H130", "246000.00"),
    @("H131", "12549.26"),
    @("H132", "2000.00"),
    @("H133", "3000.00"),
    @("H134", "1800.00"),
    @("H135", "3850.00"),
    @("H136", "500.00"),
    @("H137", "2950.00"),
    @("H138", "1600.00"),
    @("H139", "5400.00"),
    @("H140", "5665.29"),
    @("H141", "2000.00"),
    @("H142", "15705.00"),
    @("H143", "1300.00"),
    @("H144", "1127.00"),
    @("H145", "1350.00"),
    @("H146", "400.00"),
    @("H147", "17874.00"),
    @("H148", "295.80"),
    @("H149", "6340.00"),
    @("H150", "11750.00"),
    @("H151", "77.10"),
    @("H152", "2037.00"),
    @("H153", "3266.00"),
    @("H154", "16400.00"),
    @("H155", "5280.00"),
    @("H156", "5218.50"),
    @("H157", "2313.51"),
    @("H158", "5536.00"),
    @("H159", "8754.31"),
    @("H160", "7322.00"),
    @("H161", "2600.00"),
    @("H162", "6327.10"),
    @("H163", "618.00"),
    @("H164", "559.40"),
    @("H165", "1448.50"),
    @("H166", "4732.63"),
    @("H167", "11361.40"),
    @("H168", "13200.00"),
    @("H169", "1807.43"),
    @("H170", "2350.00"),
    @("H171", "2500.00"),
    @("H172", "1633.50"),
    @("H173", "117050.00"),
    @("H174", "118000.00"),
    @("H175", "262350.00"),
    @("H176", "13932.26"),
    @("H177", "26000.00"),
    @("H178", "119200.00"),
    @("H179", "2856988.11"),
    @("H180", "51364.91"),
    @("H181", "4000.00"),
    @("H182", "5790.00"),
    @("H183", "23960.00"),
    @("H184", "75780.00"),
    @("H185", "1500.00"),
    @("H186", "85330.00"),
    @("H187", "4500.00"),
    @("H188", "1653.96"),
    @("H189", "7400.00"),
    @("H190", "1000.00"),
    @("H191", "1200.00"),
)

foreach ($pair in $amountCells) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

Write-Output ("Updated {0} name cells and {1} amount cells" -f $nameCells.Length, $amountCells.Length)
